# Auto update: 2025-12-05 02:00:49
# Applies the daily refresh to DECISION/국장_방산_분석.xlsx:
#   - MACRO date label 2025-12-03 -> 2025-12-05
#   - per-ticker metrics (price, RSI, 5d return, rule score, 3/5/10d up-probabilities,
#     final score, final-score-model N column) refreshed for the 5 rows
#   - MACRO_SIGNAL judgement text switched from bullish to neutral

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: date label -------------------------------------------------
# A plain Range.Value = "2025-12-05" assignment gets auto-parsed by Excel's
# smart-typing as a date literal (because the text looks like one), which
# would turn these cells into date-formatted number cells instead of the
# plain text the workbook actually stores. Going through a TRIM() formula
# and then converting the formula to its value with PasteSpecial keeps the
# result as literal text (no reformatting / no new cell style) exactly like
# the source file.
for ($r = 2; $r -le 6; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Formula = '=TRIM("2025-12-05 ")'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = $false

# --- Row 2: KOREA AEROSPACE (047810.KS) ------------------------------------
$ws.Cells.Item(2, 4).Value  = 105500             # D  종가
$ws.Cells.Item(2, 5).Value  = 53.5               # E  RSI
$ws.Cells.Item(2, 6).Value  = -3.12              # F  5일수익률
$ws.Cells.Item(2, 7).Value  = 30                 # G  점수(룰)
$ws.Cells.Item(2, 8).Value  = 60                 # H  3일상승확률(%)
$ws.Cells.Item(2, 9).Value  = 60                 # I  5일상승확률(%)
$ws.Cells.Item(2, 10).Value = 60                 # J  10일상승확률(%)
$ws.Cells.Item(2, 11).Value = 48.7               # K  최종점수
$ws.Cells.Item(2, 14).Value = 52.43913937059539  # N  MACRO_SCORE

# --- Row 3: HYUNDAI ROTEM (064350.KS) --------------------------------------
$ws.Cells.Item(3, 4).Value  = 176200
$ws.Cells.Item(3, 5).Value  = 29.2
$ws.Cells.Item(3, 6).Value  = -1.67
$ws.Cells.Item(3, 9).Value  = 73
$ws.Cells.Item(3, 10).Value = 76
$ws.Cells.Item(3, 11).Value = 47.9
$ws.Cells.Item(3, 14).Value = 52.43913937059539

# --- Row 4: HANWHA AEROSPACE (012450.KS) -----------------------------------
$ws.Cells.Item(4, 4).Value  = 855000
$ws.Cells.Item(4, 5).Value  = 33.8
$ws.Cells.Item(4, 6).Value  = -0.93
$ws.Cells.Item(4, 8).Value  = 60
$ws.Cells.Item(4, 9).Value  = 70
$ws.Cells.Item(4, 10).Value = 73
$ws.Cells.Item(4, 11).Value = 46.7
$ws.Cells.Item(4, 14).Value = 52.43913937059539

# --- Row 5: HANWHA SYSTEMS (272210.KS) -------------------------------------
$ws.Cells.Item(5, 4).Value  = 46150
$ws.Cells.Item(5, 5).Value  = 20.9
$ws.Cells.Item(5, 6).Value  = -0.43
$ws.Cells.Item(5, 8).Value  = 56
$ws.Cells.Item(5, 9).Value  = 56
$ws.Cells.Item(5, 11).Value = 41.1
$ws.Cells.Item(5, 14).Value = 52.43913937059539

# --- Row 6: LIG Nex1 (079550.KS) --------------------------------------------
$ws.Cells.Item(6, 4).Value  = 368000
$ws.Cells.Item(6, 5).Value  = 35.8
$ws.Cells.Item(6, 6).Value  = -5.52
$ws.Cells.Item(6, 8).Value  = 66
$ws.Cells.Item(6, 10).Value = 53
$ws.Cells.Item(6, 11).Value = 39.9
$ws.Cells.Item(6, 14).Value = 52.43913937059539

# --- Column O: MACRO_SIGNAL judgement text (all 5 data rows share the text) -
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 15).Value = "⚪ 중립 구간"
}

"done"
